# Auto update ESO load data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals = @(1,30,5024,4716,4661,4666,4673,4813,5314,5880,6381,6540,6486,6428,6338,6283,6216,6205,6253,6414,6610,6589,6341,6024,5736,5298)
$col = 2
foreach ($v in $vals) {
    $ws.Cells.Item(2, $col).Value = $v
    $col = $col + 1
}

$vals = @(1,31,5142,4879,4759,4717,4740,4809,4947,5219,5782,6260,6398,6374,6305,6237,6178,6171,6253,6461,6579,6408,6124,5809,5642,5347)
$col = 2
foreach ($v in $vals) {
    $ws.Cells.Item(3, $col).Value = $v
    $col = $col + 1
}

$vals = @(2,1,5385,5152,5055,4999,4970,5053,5204,5535,6005,6469,6667,6637,6583,6504,6445,6475,6634,6918,7091,7006,6758,6416,6176,5751)
$col = 2
foreach ($v in $vals) {
    $ws.Cells.Item(4, $col).Value = $v
    $col = $col + 1
}

$vals = @(2,2,5617,5354,5232,5253,5345,5511,6021,6629,7167,7345,7243,7134,7009,6962,6922,6944,6966,7197,7494,7460,7223,6891,6589,6111)
$col = 2
foreach ($v in $vals) {
    $ws.Cells.Item(5, $col).Value = $v
    $col = $col + 1
}

$vals = @(2,3,5628,5321,5266,5272,5280,5431,5973,6586,7127,7300,7234,7162,7052,6985,6904,6891,6956,7179,7450,7428,7165,6831,6527,6065)
$col = 2
foreach ($v in $vals) {
    $ws.Cells.Item(6, $col).Value = $v
    $col = $col + 1
}

$vals = @(2,4,5400,5096,5042,5048,5056,5206,5746,6357,6896,7068,7014,6955,6865,6810,6743,6732,6790,6987,7226,7203,6939,6601,6295,5828)
$col = 2
foreach ($v in $vals) {
    $ws.Cells.Item(7, $col).Value = $v
    $col = $col + 1
}
